# Auto-generated Excel COM-interop edit script.
# Applies the "New crime data collected" weekly refresh to the CompStat sheet:
#  - bumps the report volume/number and the covered week dates
#  - refreshes every weekly/28-day/YTD/2-year crime-stat figure in rows 14-33
#  - inserts a new blank row before the old row 56 (Prepared-by / unit footer),
#    shifting it (and the trailing row) down by one, matching the new A1:N58 extent

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/number and the covered-week date range ---------------
$ws.Range("A8").Value = "Volume 31   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/16/2024  Through  12/22/2024"

# --- Numeric crime-stat refresh (rows 14-33) ------------------------------
$ws.Range("F14").Value = 1
$ws.Range("N14").Value = -64
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("G15").Value = 9
$ws.Range("H15").Value = -77.777777777777
$ws.Range("I15").Value = 62
$ws.Range("J15").Value = 43
$ws.Range("K15").Value = 44.186046511627
$ws.Range("L15").Value = 72.222222222222
$ws.Range("M15").Value = 1.639344262295
$ws.Range("N15").Value = -19.480519480519
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -56
$ws.Range("I16").Value = 288
$ws.Range("J16").Value = 311
$ws.Range("K16").Value = -7.395498392282
$ws.Range("L16").Value = 2.857142857142
$ws.Range("M16").Value = -30.769230769230
$ws.Range("N16").Value = -78.870139398385
$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 46.153846153846
$ws.Range("F17").Value = 64
$ws.Range("H17").Value = -3.030303030303
$ws.Range("I17").Value = 916
$ws.Range("J17").Value = 947
$ws.Range("K17").Value = -3.273495248152
$ws.Range("L17").Value = 17.586649550706
$ws.Range("M17").Value = 93.657505285412
$ws.Range("N17").Value = -21.508140531276
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 46
$ws.Range("H18").Value = -63.043478260869
$ws.Range("I18").Value = 335
$ws.Range("J18").Value = 363
$ws.Range("K18").Value = -7.713498622589
$ws.Range("L18").Value = -0.297619047619
$ws.Range("M18").Value = -45.439739413680
$ws.Range("N18").Value = -89.921780986762
$ws.Range("C19").Value = 28
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 97
$ws.Range("G19").Value = 106
$ws.Range("H19").Value = -8.490566037735
$ws.Range("I19").Value = 1424
$ws.Range("J19").Value = 1571
$ws.Range("K19").Value = -9.357097390197
$ws.Range("L19").Value = 3.263234227701
$ws.Range("M19").Value = 48.643006263048
$ws.Range("N19").Value = -13.170731707317
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = -78.571428571428
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 37
$ws.Range("H20").Value = -59.459459459459
$ws.Range("I20").Value = 328
$ws.Range("J20").Value = 448
$ws.Range("K20").Value = -26.785714285714
$ws.Range("L20").Value = -34.4
$ws.Range("M20").Value = -4.373177842565
$ws.Range("N20").Value = -93.117918590012
$ws.Range("C21").Value = 59
$ws.Range("D21").Value = 74
$ws.Range("E21").Value = -20.270270270270
$ws.Range("F21").Value = 207
$ws.Range("G21").Value = 289
$ws.Range("H21").Value = -28.373702422145
$ws.Range("I21").Value = 3362
$ws.Range("J21").Value = 3702
$ws.Range("K21").Value = -9.184224743381
$ws.Range("L21").Value = 1.234567901234
$ws.Range("M21").Value = 16.695591808399
$ws.Range("N21").Value = -72.803753437955
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 132
$ws.Range("K23").Value = -9.090909090909
$ws.Range("L23").Value = 27.659574468085
$ws.Range("M23").Value = 69.014084507042
$ws.Range("C24").Value = 101
$ws.Range("D24").Value = 89
$ws.Range("E24").Value = 13.483146067415
$ws.Range("F24").Value = 325
$ws.Range("G24").Value = 379
$ws.Range("H24").Value = -14.248021108179
$ws.Range("I24").Value = 3964
$ws.Range("J24").Value = 4316
$ws.Range("K24").Value = -8.155699721964
$ws.Range("L24").Value = 2.376033057851
$ws.Range("M24").Value = 7.396369547548
$ws.Range("C25").Value = 62
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = 93.75
$ws.Range("F25").Value = 184
$ws.Range("G25").Value = 130
$ws.Range("H25").Value = 41.538461538461
$ws.Range("I25").Value = 2148
$ws.Range("J25").Value = 1954
$ws.Range("K25").Value = 9.928352098259
$ws.Range("L25").Value = 51.374207188160
$ws.Range("C26").Value = 30
$ws.Range("D26").Value = 29
$ws.Range("E26").Value = 3.448275862068
$ws.Range("F26").Value = 116
$ws.Range("G26").Value = 122
$ws.Range("H26").Value = -4.918032786885
$ws.Range("I26").Value = 1830
$ws.Range("J26").Value = 1729
$ws.Range("K26").Value = 5.841526894158
$ws.Range("L26").Value = 13.453192808431
$ws.Range("M26").Value = -7.059421025901
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -81.818181818181
$ws.Range("I27").Value = 98
$ws.Range("J27").Value = 73
$ws.Range("K27").Value = 34.246575342465
$ws.Range("L27").Value = 28.947368421052
$ws.Range("C28").Value = 3
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -53.846153846153
$ws.Range("I28").Value = 194
$ws.Range("J28").Value = 208
$ws.Range("K28").Value = -6.730769230769
$ws.Range("L28").Value = 16.167664670658
$ws.Range("M29").Value = -46.428571428571
$ws.Range("N29").Value = -86.842105263157
$ws.Range("M30").Value = -42.307692307692
$ws.Range("N30").Value = -85.294117647058
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 2
$ws.Range("I31").Value = 24
$ws.Range("J31").Value = 23
$ws.Range("K31").Value = 4.347826086956
$ws.Range("L31").Value = -4
$ws.Range("L33").Value = -14.285714285714

# --- Cells that flip from a numeric 0 / computed pct to the sheet's text
# sentinels ("0" / "***.*") when the underlying count drops to nothing ------
$ws.Range("C23").Value = "'0"
$ws.Range("F29").Value = "'0"
$ws.Range("G29").Value = "'0"
$ws.Range("H29").Value = "'***.*"
$ws.Range("F30").Value = "'0"
$ws.Range("G30").Value = "'0"
$ws.Range("H30").Value = "'***.*"
$ws.Range("C31").Value = "'0"
$ws.Range("C33").Value = "'0"

# --- Insert a new blank row above the old row 56 (footer block shifts down,
# dimension grows from A1:N57 to A1:N58) ----------------------------------
$ws.Rows.Item(56).Insert()

